$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (F column) for a few events
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 2769
$ws1.Range("F9").Value = 7757
$ws1.Range("F13").Value = 329

# Sheet "全部类型" (sheet4): same events appear again, update matching rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 2769
$ws4.Range("F11").Value = 7757
$ws4.Range("F17").Value = 329
